$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 2, shifting all existing
# data rows down by one (old row N -> new row N+1).
$ws.Rows(2).Insert()

# The whole-row insert copies the header row's bold/boxed formatting onto
# the new row; clear that so the new row matches the plain data-row look,
# then restore the date number format on column D (same as every other
# data row in this sheet).
$ws.Range("A2:R2").ClearFormats()
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat

# Populate the new row 2 with the new weekly price-report record.
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Vega Modelo de Temuco"
$ws.Range("C2").Value = "La Araucanía"
$ws.Range("D2").Value = 44599
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 100114002
$ws.Range("G2").Value = "Camote"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 15
$ws.Range("K2").Value = 18000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 18000
$ws.Range("N2").Value = "$/malla 20 kilos"
$ws.Range("O2").Value = "Perú"
$ws.Range("P2").Value = 900
$ws.Range("Q2").Value = 20
$ws.Range("R2").Value = "Hortaliza"
